$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the training data values in row 1 (A1:E1) from 10 to 15
$ws.Range("A1:E1").Value = 15
